$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D9").Value = "Why SIAI – 2. 여긴 교수님들이 책 밖에 있으신 분들인거 같아서요"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-siai-2-prof-outside-ivory-tower/#utm_source=rss&utm_medium=rss&utm_campaign=why-siai-2-prof-outside-ivory-tower"

$ws.Range("D24").Value = "[근황] AAAI 2022 억셉"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222587578140"

$wb.Save()
